# Update weekly Fruta/Hortaliza price data (Repollo - Agricola del Norte S.A. de Arica)
# Columns: D=Fecha, I=Calidad, J=Volumen, K=Precio minimo, L=Precio maximo, M=Precio promedio ponderado, P=Precio $/Kg

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @{ D = 44658; I = "Segunda"; J = 1000; K = 600; L = 650; M = 625; P = 625 }
    3  = @{ D = 44474; I = "Segunda"; J = 200;  K = 600; L = 700; M = 650; P = 650 }
    4  = @{ D = 44224; I = "Segunda"; J = 800;  K = 850; L = 900; M = 875; P = 875 }
    5  = @{ D = 44201; I = "Segunda"; J = 500;  K = 800; L = 900; M = 850; P = 850 }
    6  = @{ D = 44245; I = "Primera"; J = 800;  K = 850; L = 900; M = 875; P = 875 }
    7  = @{ D = 44245; I = "Segunda"; J = 1000; K = 750; L = 800; M = 775; P = 775 }
    8  = @{ D = 44573; I = "Tercera"; J = 800;  K = 600; L = 650; M = 625; P = 625 }
    9  = @{ D = 44278; I = "Segunda"; J = 700;  K = 600; L = 700; M = 650; P = 650 }
    10 = @{ D = 44278; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 }
    11 = @{ D = 44210; I = "Segunda"; J = 900;  K = 600; L = 700; M = 650; P = 650 }
    12 = @{ D = 44544; I = "Primera"; J = 1000; K = 600; L = 650; M = 625; P = 625 }
    13 = @{ D = 44267; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 }
    14 = @{ D = 44174; I = "Segunda"; J = 800;  K = 450; L = 500; M = 475; P = 475 }
    15 = @{ D = 44174; I = "Tercera"; J = 1200; K = 250; L = 350; M = 300; P = 300 }
    16 = @{ D = 44253; I = "Segunda"; J = 1000; K = 800; L = 900; M = 850; P = 850 }
    17 = @{ D = 44253; I = "Tercera"; J = 800;  K = 600; L = 700; M = 650; P = 650 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 9).Value  = $vals.I   # I: Calidad
    $ws.Cells.Item($r, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Precio $/Kg
}
